# repositories.pptx talk deck update:
#  1. Bump the cached "datetimeFigureOut" auto-date field (footer date
#     placeholder) from 11/8/2021 to 12/8/2021 on the slide master and on
#     every slide layout that carries one.
#  2. Update the "Python 3.6+ venv" label on the GitHub-workflow picture
#     (slide 1) to "Python 3.7+ venv".

$p = $ppt.ActivePresentation

function Update-DatePlaceholders($shapes, [string]$newDate) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -like "*11/8/2021*") {
                $tr.Text = $newDate
            }
        }
    }
}

# -- 1. Slide master date placeholder --
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes "12/8/2021"

# -- 1b. Every custom (slide) layout's date placeholder --
$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DatePlaceholders $layouts.Item($j).Shapes "12/8/2021"
}

# -- 2. "Python 3.6+ venv" -> "Python 3.7+ venv" on slide 1 --
for ($n = 1; $n -le $p.Slides.Count; $n++) {
    $slide = $p.Slides.Item($n)
    for ($k = 1; $k -le $slide.Shapes.Count; $k++) {
        $shp = $slide.Shapes.Item($k)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "Python 3.6+ venv") {
                $tr.Text = "Python 3.7+ venv"
            }
        }
    }
}
